$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.714.82'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '2.351.46'
$ws.Range("E3").Value = '  -4.18%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '542.96'
$ws.Range("E5").Value = '  -1.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.59'
$ws.Range("E6").Value = '  -5.83%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.520'
$ws.Range("E8").Value = '  -11.89%  '
$ws.Range("D9").Value = '2.351.71'
$ws.Range("E9").Value = '  -3.97%  '
$ws.Range("E10").Value = '  -2.53%  '
$ws.Range("E11").Value = '  +0.16%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.23'
$ws.Range("E12").Value = '  -3.46%  '
$ws.Range("E13").Value = '  -3.55%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.69'
$ws.Range("E14").Value = '  -5.09%  '
$ws.Range("D15").Value = '2.774.50'
$ws.Range("E15").Value = '  -4.18%  '
$ws.Range("D16").Value = '60.385.29'
$ws.Range("E16").Value = '  -1.24%  '
$ws.Range("E17").Value = '  -3.18%  '
$ws.Range("D18").Value = '2.351.14'
$ws.Range("E18").Value = '  -4.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.58'
$ws.Range("E19").Value = '  -4.23%  '
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '314.03'
$ws.Range("E21").Value = '  -1.09%  '
$ws.Range("E22").Value = '  -8.37%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.87'
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.09'
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("E26").Value = '  +7.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.998'
$ws.Range("E27").Value = '  -0.33%  '
$ws.Range("D28").Value = '2.463.04'
$ws.Range("E28").Value = '  -5.20%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0891'
$ws.Range("E29").Value = '  -8.68%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.94'
$ws.Range("E30").Value = '  -3.61%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '502.07'
$ws.Range("E31").Value = '  -9.05%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("E32").Value = '  -5.96%  '
$ws.Range("E33").Value = '  -1.84%  '
$ws.Range("E34").Value = '  -5.86%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("E35").Value = '  -2.77%  '
$ws.Range("E36").Value = '  +0.07%  '
$ws.Range("E37").Value = '  -5.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.372'
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.36'
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("E40").Value = '  -10.78%  '
$ws.Range("E41").Value = '  +1.99%  '
$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.01%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '138.21'
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.09'
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.13'
$ws.Range("E45").Value = '  -10.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '138.69'
$ws.Range("E46").Value = '  -5.22%  '
$ws.Range("E47").Value = '  -1.68%  '
$ws.Range("E48").Value = '  -4.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.50'
$ws.Range("E49").Value = '  -8.59%  '
$ws.Range("E50").Value = '  -3.31%  '
$ws.Range("E51").Value = '  -4.27%  '
